# The original document has no word/styles.xml part at all (its one
# paragraph/run relies purely on direct run formatting, not a style
# reference). Touching the Styles collection forces Word to create and
# persist the styles part for the document, defining the built-in
# "Normal" paragraph style -- without altering any existing document
# content/formatting.
$d = $word.ActiveDocument

$normalStyle = $d.Styles.Add("Normal", 1)
